$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Content")
$ws2 = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------
# 1. Update the "Content" sheet data rows.
#    The validation list gained a new brand (JLoBeauty) which was
#    inserted at the top of the first brand block, pushing the
#    existing brand rows down by one. A brand-new block
#    (Prod/JLoBeauty/Core + End) was inserted further down, and the
#    trailing Smileactives/Core2 block moved from rows 20-21 to 25-26.
# ---------------------------------------------------------------

# Remove the old rows 20-21 first (their content moves to 25-26 below);
# doing this before touching anything further down keeps all later row
# numbers stable.
$ws.Range("A20:C21").EntireRow.Delete()

# Row 2: CrepeErase -> JLoBeauty
$ws.Range("B2").Value = "JLoBeauty"

# Row 3: MeaningfulBeauty -> CrepeErase
$ws.Range("B3").Value = "CrepeErase"

# Row 4: WestmoreBeauty -> MeaningfulBeauty
$ws.Range("B4").Value = "MeaningfulBeauty"

# Row 5: Smileactives/Core2 -> WestmoreBeauty/Core
$ws.Range("B5").Value = "WestmoreBeauty"
$ws.Range("C5").Value = "Core"

# Row 6: MallyBeauty/Core -> Smileactives/Core2
$ws.Range("B6").Value = "Smileactives"
$ws.Range("C6").Value = "Core2"

# Row 7: was only "End" -> now a full Prod/MallyBeauty/Core row.
# Copy the format of an existing full data row onto A7:C7 first so the
# new cells pick up the correct fill style, then set their values.
$ws.Range("A2:C2").Copy()
$ws.Range("A7:C7").PasteSpecial(-4122)
$ws.Range("A7").Value = "Prod"
$ws.Range("B7").Value = "MallyBeauty"
$ws.Range("C7").Value = "Core"

# Row 8 (new): the "End" marker that used to be on row 7.
$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "End"

# Rows 12-13 (new block): Prod/JLoBeauty/Core + End
$ws.Range("A2:C2").Copy()
$ws.Range("A12:C12").PasteSpecial(-4122)
$ws.Range("A12").Value = "Prod"
$ws.Range("B12").Value = "JLoBeauty"
$ws.Range("C12").Value = "Core"

$ws.Range("A2").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "End"

# Rows 25-26 (new position of the trailing Smileactives/Core2 block).
$ws.Range("A2:C2").Copy()
$ws.Range("A25:C25").PasteSpecial(-4122)
$ws.Range("A25").Value = "Prod"
$ws.Range("B25").Value = "Smileactives"
$ws.Range("C25").Value = "Core2"

$ws.Range("A2").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$ws.Range("A26").Value = "End"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# 2. Selections, matching the post-edit state recorded in the diff.
# ---------------------------------------------------------------
$ws2.Range("A8:XFD8").Select() | Out-Null

$ws.Activate() | Out-Null
$ws.Range("A12:XFD15").Select() | Out-Null
